$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "createing" -> "creating"
$ws.Range("B1").Value = "Test Case: Testing to see if creating assessments reaches data base with no errors"

# Step 1: add quotes around "create assessment"
$ws.Range("C2").Value = 'Step 1: While logged out go to the "create assessment" page'

# Step 2: rewrite
$ws.Range("C3").Value = "Step 2: Login as a user with the appropriate role and go to the create assessment page"

# Step 3: "a assessment" -> "an assessment"
$ws.Range("C4").Value = "Step 3: Fill out an assessment for an employee on your team"
# Expected result: "added to the database on that employee" -> "added to the database for that employee"
$ws.Range("D4").Value = "A new assessment will be added to the database for that employee"

# Step 4: "not on your team" -> "not in your team"
$ws.Range("C5").Value = "Step 4: Fill out a assessment for an employee that is not in your team"
# Expected result completely rewritten
$ws.Range("D5").Value = "I am denied access to this"

# Step 5 expected result rewritten
$ws.Range("D6").Value = "An error text shows saying that some of the fields have been left blank"

# Step 6 expected result removed -> now same text as D5 ("I am denied access to this")
$ws.Range("D7").Value = "I am denied access to this"

# Update view state: topLeftCell A2 -> A3, selection D5 -> C6
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C6").Select()
